$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.95%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'51.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'6.01%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.111"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.05%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07794"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.20%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.499"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.88%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.357"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.53%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.585"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.84%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-4.99%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1985"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.18%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'BitrueCoin"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.04736"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.40%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'MandalaExchangeToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.09510"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.01%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1045"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.10%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001260"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-5.75%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005805"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.78%"
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'2,014.50%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.12%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.440"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.02%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'1.50%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.19%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1365"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.57%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04173"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.14%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001270"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.62%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.003926"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.51%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001350"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.16%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02588"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-3.75%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05897"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'3.14%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01048"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'65.97%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008085"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.34%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1420"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.63%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008241"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.23%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008427"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'7.11%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3119"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-10.92%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007356"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.63%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05749"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'4.65%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002621"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-34.51%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
